$d = $word.ActiveDocument

# Locate the run of text "(stopping here: " that precedes the youtu.be hyperlink.
$found = $d.Content
$found.Find.Execute("(stopping here: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$runStart = $found.Start
$runEnd = $found.End

# Offsets (relative to $runStart) of the three desired run boundaries:
#   "("        -> [0, 1)
#   "stopping" -> [1, 9)
#   " here: "  -> [9, end)
$openParenEnd = $runStart + 1
$wordEnd = $runStart + 9

# First, recapitalize "stopping" -> "Stopping". Doing this via a Range.Text
# assignment re-lays-out the whole paragraph (merging same-formatted runs),
# so we do it *before* forcing the run split below.
$wordRange = $d.Range($runStart + 1, $wordEnd)
$wordRange.Text = "Stopping"

# Now force the run to split into three distinct <w:r> elements at the
# boundaries above. Toggling a character property on and back off leaves
# the visible formatting untouched but breaks the run apart from its
# neighbors, which otherwise get coalesced back together on save.
$openParen = $d.Range($runStart, $openParenEnd)
$openParen.Bold = 1
$openParen.Bold = 0

$stoppingWord = $d.Range($openParenEnd, $wordEnd)
$stoppingWord.Bold = 1
$stoppingWord.Bold = 0
